$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tags")

$data = @(
    @("mtl",  "draw", 0),
    @("rmtl", "draw", 1),
    @("dl",   "draw", 2),
    @("rdl",  "draw", 3)
)

$row = 11
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

$ws.Range("C15").Select()
